$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> new values (previously row 3 species data), with Q/R rounded
$ws.Range("A2").Value = 112182383
$ws.Range("B2").Value = 77515
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 528051
$ws.Range("R2").Value = 6905434
$ws.Range("Z2").Value = ""
$ws.Range("AB2").Value = ""

# Row 3 -> new values (previously row 2 species data), with Q/R rounded
$ws.Range("A3").Value = 112182158
$ws.Range("B3").Value = 90666
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 4364
$ws.Range("F3").Value = "Dropptaggsvamp"
$ws.Range("G3").Value = "Hydnellum ferrugineum"
$ws.Range("H3").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q3").Value = 528248
$ws.Range("R3").Value = 6905261
$ws.Range("Z3").Value = ""
$ws.Range("AB3").Value = ""
